# Economic Dashboard update - 2025-12-19
# Refresh "Latest Period"/observation dates and values for the rows whose
# source series rolled to a new reading, and move the "new data" (yellow)
# highlight to the rows that now hold the freshest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 29: 5yr, 5yr Forward (T5YIFR) ----
$ws.Range("N29").Value = 46009
$ws.Range("R29").Value = 2.22
$ws.Range("T29").Value = 2.21

# ---- Row 30: 10yr TIPS (T10YIE) ----
$ws.Range("N30").Value = 46009
$ws.Range("R30").Value = 2.24
$ws.Range("S30").Value = 2.23
$ws.Range("T30").Value = 2.25

# ---- Row 39: Nominal Broad US Dollar Index (DTWEXBGS) ----
# Date value is unchanged (46003); only the "new data" highlight moves off
# of this row, so copy the plain (non-highlighted) format from C3 onto N39.
$ws.Range("C3").Copy() | Out-Null
$ws.Range("N39").PasteSpecial(-4122) | Out-Null

# ---- Row 47: FFR (DFF) ----
$ws.Range("N47").Value = 46008

# ---- Row 48: 2y UST (DGS2) ----
$ws.Range("N48").Value = 46008
$ws.Range("Q48").Value = 3.49
$ws.Range("R48").Value = 3.48
$ws.Range("S48").Value = 3.51
$ws.Range("U48").Value = ""

# ---- Row 49: 5y UST (DGS5) ----
$ws.Range("N49").Value = 46008
$ws.Range("Q49").Value = 3.7
$ws.Range("R49").Value = 3.69
$ws.Range("S49").Value = 3.73
$ws.Range("U49").Value = ""

# ---- Row 50: 10y UST (DGS10) ----
$ws.Range("N50").Value = 46008
$ws.Range("Q50").Value = 4.16
$ws.Range("R50").Value = 4.15
$ws.Range("S50").Value = 4.18
$ws.Range("U50").Value = ""

# ---- Row 51: 30y Mtg. (MORTGAGE30US) ----
# This row now carries the freshest observation, so it picks up the
# "new data" highlight - copy the highlighted format from N47/N48 onto N51.
$ws.Range("N47").Copy() | Out-Null
$ws.Range("N51").PasteSpecial(-4122) | Out-Null
$ws.Range("N51").Value = 46006
$ws.Range("Q51").Value = 6.21
$ws.Range("R51").Value = 6.22
$ws.Range("S51").Value = 6.19
$ws.Range("T51").Value = 6.23
$ws.Range("U51").Value = 6.26

# ---- Row 52: BAA (DBAA) ----
$ws.Range("N52").Value = 46008
$ws.Range("Q52").Value = 5.94
$ws.Range("R52").Value = 5.93
$ws.Range("S52").Value = 5.95
$ws.Range("U52").Value = ""
